$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the N column (SPT blow count) values for rows 11-25
$ws.Range("N11").Value = 60
$ws.Range("N12").Value = 60
$ws.Range("N13").Value = 80
$ws.Range("N14").Value = 100
$ws.Range("N15").Value = 100
$ws.Range("N16").Value = 100
$ws.Range("N17").Value = 60
$ws.Range("N18").Value = 60
$ws.Range("N19").Value = 60
$ws.Range("N20").Value = 90
$ws.Range("N21").Value = 98
$ws.Range("N22").Value = 98
$ws.Range("N23").Value = 99
$ws.Range("N24").Value = 99
$ws.Range("N25").Value = 100

# Update the active selection to N18 (matches saved UI state in the diff)
$ws.Range("N18").Select()
